$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 513, shifting existing rows 513:579 down to 514:580.
$ws.Rows.Item(513).Insert()

# Populate the newly inserted row 513 with its full record, matching the
# surrounding rows' constant columns and the new Fecha/Volumen values.
$ws.Cells.Item(513, 1).Value = 10
$ws.Cells.Item(513, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(513, 3).Value = "La Araucanía"
$ws.Cells.Item(513, 4).Value = 45212
$ws.Cells.Item(513, 5).Value = 9
$ws.Cells.Item(513, 6).Value = 100112009
$ws.Cells.Item(513, 7).Value = "Acelga"
$ws.Cells.Item(513, 8).Value = "Sin especificar"
$ws.Cells.Item(513, 9).Value = "Primera"
$ws.Cells.Item(513, 10).Value = 40
$ws.Cells.Item(513, 11).Value = 8000
$ws.Cells.Item(513, 12).Value = 8000
$ws.Cells.Item(513, 13).Value = 8000
$ws.Cells.Item(513, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(513, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(513, 16).Value = 667
$ws.Cells.Item(513, 17).Value = 12
$ws.Cells.Item(513, 18).Value = "Hortaliza"
